# Continuing Issue724 - tweak DevMan class/package diagram shapes and
# register an ARTICULATE_PROJECT_OPEN custom tag on the presentation.
#
# NOTE: Shape.Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU)
# in the PowerPoint object model, while the underlying OOXML stores English
# Metric Units (EMU). The literal point values below were chosen so that,
# after the host's internal float32 round-trip, they serialize back to the
# exact target EMU values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Presentation-level custom tag (adds p:custDataLst / ppt/tags/tag1.xml) ---
$p.Tags.Add("ARTICULATE_PROJECT_OPEN", "0")

# --- Shape 1: "Rounded Rectangle 88" (id=89) - move/resize left edge ---
$sh1 = $s.Shapes.Item(1)
$sh1.Left = 434181 / 12700
$sh1.Width = 1981200 / 12700

# --- Shape 15: "Straight Arrow Connector 114" (id=115) - reposition ---
$sh15 = $s.Shapes.Item(15)
$sh15.Left = 88.18748094488188
$sh15.Top = 914402 / 12700
$sh15.Width = 0
$sh15.Height = 475534 / 12700

# --- Shape 18: "Straight Arrow Connector 117" (id=118) - shift right ---
$sh18 = $s.Shapes.Item(18)
$sh18.Left = 662781 / 12700

# --- Shape 38: "Group 85" (id=86) - move up slightly ---
$sh38 = $s.Shapes.Item(38)
$sh38.Top = 5220494 / 12700

# Inside that group, retype "client::scripts" to drop the stray trailing
# endParaRPr while keeping the run's formatting (color/bold/size) intact.
$rect86 = $sh38.GroupItems.Item(1)
$rect86.TextFrame.TextRange.Delete()
$rect86.TextFrame.TextRange.Text = "client::scripts"

# --- Shape 39: "Straight Arrow Connector 161" bent connector (id=162) ---
$sh39 = $s.Shapes.Item(39)
$sh39.Height = 877094 / 12700

# --- Shape 40: "Group 168" (id=169) - move down slightly ---
$sh40 = $s.Shapes.Item(40)
$sh40.Top = 5779294 / 12700

# --- Shape 42: "Straight Arrow Connector 161" (id=176) - shift right ---
$sh42 = $s.Shapes.Item(42)
$sh42.Left = 281781 / 12700

# --- New connector: duplicate shape 17 ("Straight Arrow Connector 116"),
# which already has the same un-flipped line style/formatting we need
# (00B050 sysDot line, none/triangle arrowheads), then reposition/rename.
# Duplicate() appends the copy at the end of the shape collection, which
# lands it right after shape 176 - matching the target insertion point.
$dup = $s.Shapes.Item(17).Duplicate()
$newConn = $dup.Item(1)
$newConn.Name = "Straight Arrow Connector 77"
$newConn.Left = 76.1874809448819
$newConn.Top = 443.23686220472445
$newConn.Width = 0
$newConn.Height = 17.651339606299214
